# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45883

$ws.Range("B2").Value = 126.87
$ws.Range("C2").Value = 115.7
$ws.Range("D2").Value = 110.28
$ws.Range("E2").Value = 105.27
$ws.Range("F2").Value = 102
$ws.Range("G2").Value = 101.12
$ws.Range("H2").Value = 107
$ws.Range("I2").Value = 109.65
$ws.Range("J2").Value = 99.31
$ws.Range("K2").Value = 90.92
$ws.Range("L2").Value = 71.67
$ws.Range("M2").Value = 67.40000000000001
$ws.Range("N2").Value = 54.01
$ws.Range("O2").Value = 40.91
$ws.Range("P2").Value = 54.9
$ws.Range("Q2").Value = 69.43000000000001
$ws.Range("R2").Value = 71.67
$ws.Range("S2").Value = 80
$ws.Range("T2").Value = 102
$ws.Range("U2").Value = 115.78
$ws.Range("V2").Value = 137.99
$ws.Range("W2").Value = 149.02
$ws.Range("X2").Value = 125.88
$ws.Range("Y2").Value = 111.92
$ws.Range("Z2").Value = 96.7

$ws.Range("AB2").Value = 131.2
$ws.Range("AD2").Value = 143.5
$ws.Range("AE2").Value = "0h-2h"
$ws.Range("AF2").Value = 121.28
